$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values are stored as text (inlineStr) in the workbook, so force the
# cells to text format first to avoid Excel auto-converting numeric-looking
# strings into real numbers.
$textCells = @("C2", "E2", "F2", "I2", "P2", "Q2")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("C2").Value = "8761"
$ws.Range("E2").Value = "0.224"
$ws.Range("F2").Value = "1.409"
$ws.Range("I2").Value = "1693"
$ws.Range("P2").Value = "95.219"
$ws.Range("Q2").Value = "49757.45"
